# Add a new row 5 to Sheet1 with a "TEST" calendar entry, mirroring the
# structure of the existing row 4 (TEST/EVENT row), and move the active
# selection to N5 (one column past the last used column on the new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns
$ws.Range("A5").Value = "FINAL TEST"
$ws.Range("B5").Value = "TEST"
$ws.Range("C5").Value = "TEST"
$ws.Range("D5").Value = "TEST"

# Start date / start time / end date (same pattern as row 4: date, blank
# time cell, date) reusing the workbook's existing date/time number formats.
$ws.Range("E5").Value = 45884
$ws.Range("E5").NumberFormat = "d-mmm"

$ws.Range("F5").NumberFormat = "h:mm"

$ws.Range("G5").Value = 45885
$ws.Range("G5").NumberFormat = "d-mmm"

# Timezone columns + link/description column
$ws.Range("I5").Value = "Melbourne"
$ws.Range("J5").Value = "Melbourne"
$ws.Range("M5").Value = "Opaque"

# Move the selection to N5, like in the authored workbook.
$null = $ws.Range("N5").Select()
